$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.290.71"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.944.08"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.78%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("E11").Value = "  +6.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "4.563.54"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "3.964.98"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "67.553.36"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("E22").Value = "  +3.12%  "
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.83%  "
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "728.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.152"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "0.0₃0783"
$ws.Range("E38").Value = "  +12.83%  "
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0477"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.52%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.336"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.06%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("E47").Value = "  +4.75%  "
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "147.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("E51").Value = "  +1.33%  "
